# Generate Report for Handoff
# Rotate the localization-status report from the old source file
# (6580e3be-e665-4aa3-ac99-b6a47b44c08f.md) to the newly generated one
# (e7a4e73a-68d0-48e9-a82c-b6de2683f2bf.md), and refresh the handoff
# timestamps / target-file names that go with it.

$wb = $excel.ActiveWorkbook

$oldGuid = "6580e3be-e665-4aa3-ac99-b6a47b44c08f"
$newGuid = "e7a4e73a-68d0-48e9-a82c-b6de2683f2bf"

$newMdName   = $newGuid + ".md"
$newMdPath   = "e2e\" + $newGuid + ".md"
$newZhXlf    = $newGuid + ".324a880f8954ffbdc76951ba4765a752eb6a73a0.zh-cn.xlf"
$newDeXlf    = $newGuid + ".324a880f8954ffbdc76951ba4765a752eb6a73a0.de-de.xlf"

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = "2016-08-28 08:56:33"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newMdPath
}

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-28 08:56:28"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-08-28 08:56:33"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}
